$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 92
$ws.Cells.Item(92,1).Value = 1
$ws.Cells.Item(92,2).Value = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(92,3).Value = 'Arica y Parinacota'
$ws.Cells.Item(92,4).Value = 45147
$ws.Cells.Item(92,5).Value = 15
$ws.Cells.Item(92,6).Value = 'Fruta'
$ws.Cells.Item(92,7).Value = 100101
$ws.Cells.Item(92,8).Value = 'Berries'
$ws.Cells.Item(92,9).Value = 100112025
$ws.Cells.Item(92,10).Value = 'Frutilla'
$ws.Cells.Item(92,11).Value = 'Sin especificar'
$ws.Cells.Item(92,12).Value = 'Primera'
$ws.Cells.Item(92,13).Value = 100
$ws.Cells.Item(92,14).Value = 6000
$ws.Cells.Item(92,15).Value = 7000
$ws.Cells.Item(92,16).Value = 6500
$ws.Cells.Item(92,17).Value = '$/bandeja 3 kilos'
$ws.Cells.Item(92,18).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(92,19).Value = 2167
$ws.Cells.Item(92,20).Value = 3
$ws.Cells.Item(92,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 93
$ws.Cells.Item(93,1).Value = 1
$ws.Cells.Item(93,2).Value = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(93,3).Value = 'Arica y Parinacota'
$ws.Cells.Item(93,4).Value = 45147
$ws.Cells.Item(93,5).Value = 15
$ws.Cells.Item(93,6).Value = 'Fruta'
$ws.Cells.Item(93,7).Value = 100101
$ws.Cells.Item(93,8).Value = 'Berries'
$ws.Cells.Item(93,9).Value = 100112025
$ws.Cells.Item(93,10).Value = 'Frutilla'
$ws.Cells.Item(93,11).Value = 'Sin especificar'
$ws.Cells.Item(93,12).Value = 'Segunda'
$ws.Cells.Item(93,13).Value = 160
$ws.Cells.Item(93,14).Value = 4000
$ws.Cells.Item(93,15).Value = 5000
$ws.Cells.Item(93,16).Value = 4500
$ws.Cells.Item(93,17).Value = '$/bandeja 3 kilos'
$ws.Cells.Item(93,18).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(93,19).Value = 1500
$ws.Cells.Item(93,20).Value = 3
$ws.Cells.Item(93,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 94
$ws.Cells.Item(94,1).Value = 1
$ws.Cells.Item(94,2).Value = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(94,3).Value = 'Arica y Parinacota'
$ws.Cells.Item(94,4).Value = 45147
$ws.Cells.Item(94,5).Value = 15
$ws.Cells.Item(94,6).Value = 'Fruta'
$ws.Cells.Item(94,7).Value = 100101
$ws.Cells.Item(94,8).Value = 'Berries'
$ws.Cells.Item(94,9).Value = 100112025
$ws.Cells.Item(94,10).Value = 'Frutilla'
$ws.Cells.Item(94,11).Value = 'Sin especificar'
$ws.Cells.Item(94,12).Value = 'Tercera'
$ws.Cells.Item(94,13).Value = 200
$ws.Cells.Item(94,14).Value = 2000
$ws.Cells.Item(94,15).Value = 3000
$ws.Cells.Item(94,16).Value = 2500
$ws.Cells.Item(94,17).Value = '$/bandeja 3 kilos'
$ws.Cells.Item(94,18).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(94,19).Value = 833
$ws.Cells.Item(94,20).Value = 3
$ws.Cells.Item(94,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
